$d = $word.ActiveDocument

# 1. "FROM" signer name: Vu Thien Nhuong -> Vu Tuan Khanh
$d.Content.Find.Execute("Vũ Thiện Nhượng", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Vũ Tuấn Khanh", 2)

# 2. Remove the ", PTGĐ. Vũ Tuấn Khanh" part after "Nguyễn Văn Nam"
$d.Content.Find.Execute("Nguyễn Văn Nam, PTGĐ. Vũ Tuấn Khanh", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Nguyễn Văn Nam", 2)

# 3. Placeholder field rename: <ContractSiteId> -> <SiteName>
$d.Content.Find.Execute("<ContractSiteId>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<SiteName>", 2)

# 4. Signature block: THIỆN NHƯỢNG -> TUẤN KHANH
$d.Content.Find.Execute("THIỆN NHƯỢNG", $true, $false, $false, $false, $false,
                         $true, 1, $false, "TUẤN KHANH", 2)
